# Auto-generated: applies numeric profit-recalculation updates to Gungnir_Profits workbook
# Sheets affected: ALC, ARM, CRP, CUL, GSM, LTW
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2015729.8
$ws.Range("I40").Value = 2606032.8
$ws.Range("J40").Value = 835123.7
$ws.Range("K40").Value = 2606032.8
$ws.Range("L40").Value = 835123.7
$ws.Range("M40").Value = -2605857.8
$ws.Range("N40").Value = -835473.7
$ws.Range("H64").Value = 3710
$ws.Range("I64").Value = 3516.6667
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 3516.6667
$ws.Range("L64").Value = 4000
$ws.Range("M64").Value = -3268.6667
$ws.Range("N64").Value = -4496
$ws.Range("H67").Value = 3710
$ws.Range("I67").Value = 3516.6667
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 3516.6667
$ws.Range("L67").Value = 4000
$ws.Range("M67").Value = -2658.6667
$ws.Range("N67").Value = -5716
$ws.Range("H74").Value = 3511.5557
$ws.Range("I74").Value = 3800
$ws.Range("J74").Value = 3429.1428
$ws.Range("K74").Value = 3800
$ws.Range("L74").Value = 3429.1428
$ws.Range("M74").Value = -2864
$ws.Range("N74").Value = -5301.1428
$ws.Range("H76").Value = 9265059
$ws.Range("I76").Value = 8895.294
$ws.Range("J76").Value = 17546890
$ws.Range("K76").Value = 8895.294
$ws.Range("L76").Value = 17546890
$ws.Range("M76").Value = -8580.294
$ws.Range("N76").Value = -17547520
$ws.Range("H77").Value = 3511.5557
$ws.Range("I77").Value = 3800
$ws.Range("J77").Value = 3429.1428
$ws.Range("K77").Value = 19000
$ws.Range("L77").Value = 17145.714
$ws.Range("M77").Value = -14320
$ws.Range("N77").Value = -26505.714
$ws.Range("H79").Value = 9265059
$ws.Range("I79").Value = 8895.294
$ws.Range("J79").Value = 17546890
$ws.Range("K79").Value = 8895.294
$ws.Range("L79").Value = 17546890
$ws.Range("M79").Value = -7803.294
$ws.Range("N79").Value = -17549074

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6945834.5
$ws.Range("I61").Value = 7937925
$ws.Range("K61").Value = 7937925
$ws.Range("M61").Value = -7937713
$ws.Range("H74").Value = 1367.0322
$ws.Range("I74").Value = 1443.75
$ws.Range("J74").Value = 1227.5454
$ws.Range("K74").Value = 1443.75
$ws.Range("L74").Value = 1227.5454
$ws.Range("M74").Value = -569.75
$ws.Range("N74").Value = -2975.5454
$ws.Range("H77").Value = 1367.0322
$ws.Range("I77").Value = 1443.75
$ws.Range("J77").Value = 1227.5454
$ws.Range("K77").Value = 7218.75
$ws.Range("L77").Value = 6137.727
$ws.Range("M77").Value = -2850.75
$ws.Range("N77").Value = -14873.727
$ws.Range("H136").Value = 6945834.5
$ws.Range("I136").Value = 7937925
$ws.Range("K136").Value = 23813775
$ws.Range("M136").Value = -23811225

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4347.1113
$ws.Range("I62").Value = 2879.8
$ws.Range("J62").Value = 6181.25
$ws.Range("K62").Value = 2879.8
$ws.Range("L62").Value = 6181.25
$ws.Range("M62").Value = -2255.8
$ws.Range("N62").Value = -7429.25
$ws.Range("H65").Value = 4347.1113
$ws.Range("I65").Value = 2879.8
$ws.Range("J65").Value = 6181.25
$ws.Range("K65").Value = 14399
$ws.Range("L65").Value = 30906.25
$ws.Range("M65").Value = -11279
$ws.Range("N65").Value = -37146.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 235715.28
$ws.Range("J37").Value = 235715.28
$ws.Range("L37").Value = 707145.84
$ws.Range("N37").Value = -707369.84
$ws.Range("H55").Value = 3800
$ws.Range("J55").Value = 3800
$ws.Range("L55").Value = 11400
$ws.Range("N55").Value = -11754
$ws.Range("H68").Value = 6190.3335
$ws.Range("I68").Value = 428.125
$ws.Range("J68").Value = 10800.1
$ws.Range("K68").Value = 1284.375
$ws.Range("L68").Value = 32400.3
$ws.Range("M68").Value = -473.375
$ws.Range("N68").Value = -34022.3
$ws.Range("H71").Value = 6190.3335
$ws.Range("I71").Value = 428.125
$ws.Range("J71").Value = 10800.1
$ws.Range("K71").Value = 3853.125
$ws.Range("L71").Value = 97200.90000000001
$ws.Range("M71").Value = 202.875
$ws.Range("N71").Value = -105312.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18829.715
$ws.Range("I70").Value = 18829.715
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 18829.715
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -18559.715
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 18829.715
$ws.Range("I73").Value = 18829.715
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 18829.715
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -17893.715
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 7694646
$ws.Range("I80").Value = 2530
$ws.Range("J80").Value = 33335034
$ws.Range("K80").Value = 2530
$ws.Range("L80").Value = 33335034
$ws.Range("M80").Value = -1532
$ws.Range("N80").Value = -33337030
$ws.Range("H83").Value = 7694646
$ws.Range("I83").Value = 2530
$ws.Range("J83").Value = 33335034
$ws.Range("K83").Value = 12650
$ws.Range("L83").Value = 166675170
$ws.Range("M83").Value = -7658
$ws.Range("N83").Value = -166685154

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 7667.75
$ws.Range("J62").Value = 5220
$ws.Range("L62").Value = 5220
$ws.Range("N62").Value = -6468
$ws.Range("H65").Value = 7667.75
$ws.Range("J65").Value = 5220
$ws.Range("L65").Value = 15660
$ws.Range("N65").Value = -21900
$ws.Range("H74").Value = 25260
$ws.Range("I74").Value = 40000
$ws.Range("K74").Value = 40000
$ws.Range("M74").Value = -39002
$ws.Range("H75").Value = 33986.5
$ws.Range("J75").Value = 33986.5
$ws.Range("L75").Value = 33986.5
$ws.Range("N75").Value = -35858.5
$ws.Range("H77").Value = 25260
$ws.Range("I77").Value = 40000
$ws.Range("K77").Value = 120000
$ws.Range("M77").Value = -115008
$ws.Range("H78").Value = 33986.5
$ws.Range("J78").Value = 33986.5
$ws.Range("L78").Value = 101959.5
$ws.Range("N78").Value = -111319.5
